# Updates cryptos list data: price (D) and 1h volume change (E) columns,
# plus a row reorder/content swap for rows 36-37 (ImmutableX / LidoDAOToken).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.711.09"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").Value = "1.600.93"
$ws.Range("E3").Value = "  +0.30%  "
$ws.Range("E4").Value = "  +0.35%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.41"
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("E6").Value = "  -0.74%  "
$ws.Range("E7").Value = "  +0.31%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.63"
$ws.Range("E10").Value = "  +0.69%  "
$ws.Range("E11").Value = "  +0.58%  "
$ws.Range("D12").Value = "1.826.00"
$ws.Range("E12").Value = "  +0.31%  "
$ws.Range("D13").Value = "1.602.86"
$ws.Range("E13").Value = "  +0.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.04"
$ws.Range("E14").Value = "  +0.28%  "
$ws.Range("E15").Value = "  +0.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.98"
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("D17").Value = "26.685.37"
$ws.Range("E17").Value = "  +0.34%  "
$ws.Range("D18").Value = "0.0₃0744"
$ws.Range("E18").Value = "  +0.54%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "210.09"
$ws.Range("E19").Value = "  +0.85%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.19"
$ws.Range("E20").Value = "  +2.93%  "
$ws.Range("E21").Value = "  +0.31%  "
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("E23").Value = "  -0.39%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.97"
$ws.Range("E24").Value = "  +0.78%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.22"
$ws.Range("E25").Value = "  -0.82%  "
$ws.Range("E26").Value = "  +0.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.11"
$ws.Range("E28").Value = "  -0.94%  "
$ws.Range("E29").Value = "  +0.55%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0513"
$ws.Range("E30").Value = "  +0.40%  "
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("E32").Value = "  +1.27%  "
$ws.Range("E33").Value = "  +0.80%  "
$ws.Range("D34").Value = "1.294.70"
$ws.Range("E34").Value = "  +0.94%  "
$ws.Range("E35").Value = "  +0.63%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.605"
$ws.Range("E36").Value = "  -2.47%  "
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.50"
$ws.Range("E37").Value = "  +0.68%  "
$ws.Range("E38").Value = "  +15.04%  "
$ws.Range("E39").Value = "  -0.51%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.824"
$ws.Range("E40").Value = "  -1.80%  "
$ws.Range("E41").Value = "  -1.36%  "
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.03"
$ws.Range("E44").Value = "  -1.88%  "
$ws.Range("D45").Value = "1.738.85"
$ws.Range("E45").Value = "  +0.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "90.59"
$ws.Range("E46").Value = "  +0.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.55"
$ws.Range("E47").Value = "  -3.03%  "
$ws.Range("E48").Value = "  -1.18%  "
$ws.Range("E49").Value = "  +1.94%  "
$ws.Range("E50").Value = "  +0.08%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.43"
$ws.Range("E51").Value = "  -0.05%  "
